$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (Price) cells are treated as text so values like "590.42" or "0.0747"
# are not auto-converted to numbers by Excel, matching the source data which stores
# prices as plain text strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.633.01'
$ws.Range("E2").Value = '  +4.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.494.25'
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.42'
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.74'
$ws.Range("E6").Value = '  +4.40%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  +8.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.492.64'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  +7.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.33'
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("E12").Value = '  +4.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.100.45'
$ws.Range("E13").Value = '  +2.69%  '
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.21'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("E16").Value = '  +3.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.642.79'
$ws.Range("E17").Value = '  +4.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.506.98'
$ws.Range("E18").Value = '  +3.34%  '
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.15'
$ws.Range("E20").Value = '  +4.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.31'
$ws.Range("E21").Value = '  +4.53%  '
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.43'
$ws.Range("E23").Value = '  +3.57%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.537'
$ws.Range("E25").Value = '  +4.56%  '
$ws.Range("E26").Value = '  +5.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.17'
$ws.Range("E27").Value = '  +7.48%  '
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +4.64%  '
$ws.Range("E31").Value = '  +6.71%  '
$ws.Range("E32").Value = '  +3.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.61'
$ws.Range("E33").Value = '  +3.40%  '
$ws.Range("E34").Value = '  +5.38%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +9.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.55'
$ws.Range("E37").Value = '  +1.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.883'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  +6.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.70'
$ws.Range("E40").Value = '  +7.42%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0747'
$ws.Range("E41").Value = '  +3.18%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.48'
$ws.Range("E42").Value = '  +5.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.73'
$ws.Range("E43").Value = '  +4.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.54'
$ws.Range("E44").Value = '  +3.14%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.20'
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.782.79'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0313'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '352.62'
$ws.Range("E48").Value = '  +7.19%  '
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("E50").Value = '  +5.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.00'
$ws.Range("E51").Value = '  +14.27%  '
